$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# player/team name + runs/balls/fours/sixes per row (row 13 is a newly added match)
$data = @(
  @(2, "MS Dhoni (c)†", "Chennai Super Kings", "16", "16", "2", "1"),
  @(3, "MS Dhoni (c)†", "Chennai Super Kings", "1", "4", "0", "0"),
  @(4, "MS Dhoni (c)†", "Chennai Super Kings", "19", "21", "3", "0"),
  @(5, "MS Dhoni (c)†", "Chennai Super Kings", "11", "12", "1", "0"),
  @(6, "MS Dhoni (c)†", "Chennai Super Kings", "29", "17", "0", "3"),
  @(7, "MS Dhoni (c)†", "Chennai Super Kings", "28", "28", "2", "0"),
  @(8, "MS Dhoni (c)†", "Chennai Super Kings", "21", "13", "2", "1"),
  @(9, "MS Dhoni (c)†", "Chennai Super Kings", "3", "5", "0", "0"),
  @(10, "MS Dhoni (c)†", "Chennai Super Kings", "47", "36", "4", "1"),
  @(11, "MS Dhoni (c)†", "Chennai Super Kings", "0", "2", "0", "0"),
  @(12, "MS Dhoni (c)†", "Chennai Super Kings", "10", "6", "0", "1"),
  @(13, "MS Dhoni (c)†", "Chennai Super Kings", "15", "12", "2", "0")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = "'" + $row[3]
    $ws.Cells.Item($r, 4).Value = "'" + $row[4]
    $ws.Cells.Item($r, 5).Value = "'" + $row[5]
    $ws.Cells.Item($r, 6).Value = "'" + $row[6]
}

